$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.375.09"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "2.120.26"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "350.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.012"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5235"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4519"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.17"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09055"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.185"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.76%  "
$ws.Range("D13").Value = "2.132.93"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.815"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.101"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001161"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.012"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06729"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.009"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.293"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "30.438.95"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.371"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "2.385.17"
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.573"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.195"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1075"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.683"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.336"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.049"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.047"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02622"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06911"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2334"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6914"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.266"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6470"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.305"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000368"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.702"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.240"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07266"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.23%  "
